$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Employee name ---
$ws.Range("C2").Value = "Priyanka Muddana"

# --- DOJ date: must stay a literal text string "2014-03-03", not be
# auto-converted into a date serial number by Excel's smart-entry parsing.
# Trick: force the cell to Text format before assigning, then restore the
# original (General, style index 0) formatting by copying it over from a
# neighboring cell that already uses that style.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2014-03-03"
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats

# --- Payroll figures (row 2) ---
$ws.Range("H2").Value = 20000.0
$ws.Range("I2").Value = 1666.6666666666667
$ws.Range("J2").Value = 1557.51
$ws.Range("M2").Value = 666.667
$ws.Range("N2").Value = 166.67
$ws.Range("O2").Value = 724.17
$ws.Range("Q2").Value = 1557.51
$ws.Range("R2").Value = 80.0
$ws.Range("S2").Value = 27.26
$ws.Range("W2").Value = 107.26
$ws.Range("X2").Value = 1450.25
